$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new "team record" columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style already used by the other header cells (e.g. AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill every data row (2-48) with the team's 2005 record: 74 wins, 88 losses, 0 ties
for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 30).Value = 74
    $ws.Cells.Item($row, 31).Value = 88
    $ws.Cells.Item($row, 32).Value = 0
}
